$d = $word.ActiveDocument

# 1. Replace the paragraph text in the table cell (problem 4 solution).
$d.Content.Find.Execute(
    "A discrete random variable is something that varies following a specific pattern",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A probability of 1 implies an event is certain to happen. A probability of 0",
    2)

$d.Content.Find.Execute(
    "or distribution over the long run. They are discrete if they can be listed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "implies it is impossible to happen, or certain to not happen.",
    2)

# 2. Adjust the table grid column widths (Width is in points; XML stores
#    twips, i.e. 1/20 pt, so divide the target twip values by 20).
$table = $d.Tables.Item(1)
$table.Columns.Item(1).Width = 800 / 20
$table.Columns.Item(2).Width = 560 / 20
$table.Columns.Item(3).Width = 6560 / 20
